$d = $word.ActiveDocument

# --- Locate the question paragraph ("Jika kedua kelas diatas ...") ---
$count = $d.Paragraphs.Count
$questionIndex = 0
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -match "Jika kedua kelas diatas") {
        $questionIndex = $i
    }
}

# The paragraph right after the question is already present (empty, with
# the formatting we want to reuse: ListParagraph style, spacing, 1080
# indent, both-justify, bold paragraph mark but no numbering). Insert a
# brand-new paragraph right before it and it will inherit that exact
# formatting.
$followingPara = $d.Paragraphs.Item($questionIndex + 1)
$followingPara.Range.InsertParagraphBefore()

$answerPara = $d.Paragraphs.Item($questionIndex + 1)
$answerRange = $answerPara.Range

# Type the answer text, plus a trailing placeholder character so the
# insertion point used for the bookmark below is never the very last
# slot of the paragraph (that boundary position is unreliable for
# Bookmarks.Add in this runtime).
$answerText = "Jawab: Tetap berfungsi, karena walaupun di package yang sama namun jika suatu class menggunakan hak akses private, maka atribut/method tersebut hanya bisa diakses dalam classnya sendiri."
$answerRange.InsertBefore($answerText + "#")

# The run itself should not be bold (only the paragraph mark keeps the
# inherited bold), so unbold everything we just typed (text + placeholder).
$typedRange = $d.Range($answerRange.Start, $answerRange.End - 1)
$typedRange.Font.Bold = 0

# Drop the bookmark right before the placeholder character.
$placeholderPos = $answerRange.End - 2
$bookmarkRange = $d.Range($placeholderPos, $placeholderPos)

# Remove the old _GoBack bookmark (was sitting alone in the previous
# empty paragraph) before re-adding it here, so there is only ever one.
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Remove the placeholder character; the now-collapsed bookmark stays put
# right after the answer text and before the paragraph mark.
$d.Range($placeholderPos, $placeholderPos + 1).Delete()
